$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Fix" the conciliated highlight color: yellow -> light green (D8E4BC)
#    for the rows that are still flagged as "conciliada" (rows 2-6, ids 1-5).
$newColor = 216 + (228 * 256) + (188 * 65536)   # RGB(216,228,188) = D8E4BC
$ws.Range("A2:K6").Interior.Color = $newColor

# 2. Rows 7-17 (ids 6-16) are no longer "conciliada":
#    - drop the highlight formatting entirely (back to default / no fill)
#    - set column K ("conciliada") to 0
$ws.Range("A7:K17").ClearFormats()
$ws.Range("K7:K17").Value = 0

# 3. Once the style that was keeping the genuinely-blank cells "alive" is
#    gone, those blank cells need to disappear completely instead of
#    lingering as empty, unstyled <c> entries.
$blankCells = @("D7","H7","D8","G8","D9","G9","D10","G10","D11","G11","D12","G12","D13","H13","D14","H14","D15","H15","D16","H16","D17","H17")
foreach ($addr in $blankCells) {
    $ws.Range($addr).Clear()
}
